# Add Indian MF 1st Stab - insert 9 new weekly-date columns at the front
# of the ratings history table, fill in their header labels, and mark the
# newly-learned rating actions for the affected research firms.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 9 new blank columns before column B (shifts B:V -> K:AE) ---
$ws.Range("B:J").Insert()

# --- 2. New weekly date headers in row 1 (most-recent-first ordering) ---
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# --- 3. Fill the new columns (B:J) for every data row with the same
#        "UN" filler used throughout the rest of the table ---
$ws.Range("B2:J33").Value = "UN"

# --- 4. New rating-change events recorded for three firms ---
# Row 5 - Zacks Investment Research
$ws.Range("B5").Value = "9/7/2019,Upgrades,Sell -> Hold,"
$ws.Range("B5").Interior.Color = 13434828
$ws.Range("C5").Value = "8/22/2019,Upgrades,Sell -> Hold,"
$ws.Range("C5").Interior.Color = 13434828

# Row 11 - ValuEngine
$ws.Range("B11").Value = "9/5/2019,Upgrades,Hold -> Buy,"
$ws.Range("B11").Interior.Color = 13434828

# Row 22 - BidaskClub
$ws.Range("B22").Value = "9/6/2019,Upgrades,Hold -> Buy,"
$ws.Range("B22").Interior.Color = 13434828
